$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.249.12"
$ws.Range("E2").Value = "  +1.17%  "
$ws.Range("D3").Value = "2.457.15"
$ws.Range("E3").Value = "  +1.25%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'574.05"
$ws.Range("E5").Value = "  +1.43%  "
$ws.Range("D6").Value = "'146.60"
$ws.Range("E6").Value = "  +1.03%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  +1.25%  "
$ws.Range("D9").Value = "2.454.09"
$ws.Range("E9").Value = "  +1.18%  "
$ws.Range("E10").Value = "  +1.43%  "
$ws.Range("D11").Value = "'0.157"
$ws.Range("E11").Value = "  +1.47%  "
$ws.Range("E12").Value = "  -0.43%  "
$ws.Range("D13").Value = "'0.357"
$ws.Range("E13").Value = "  +0.46%  "
$ws.Range("D14").Value = "'27.09"
$ws.Range("E14").Value = "  +1.01%  "
$ws.Range("E15").Value = "  +0.51%  "
$ws.Range("D16").Value = "2.901.46"
$ws.Range("E16").Value = "  +1.26%  "
$ws.Range("D17").Value = "63.089.08"
$ws.Range("E17").Value = "  +1.21%  "
$ws.Range("D18").Value = "2.452.31"
$ws.Range("E18").Value = "  +1.28%  "
$ws.Range("D19").Value = "'11.33"
$ws.Range("E19").Value = "  +1.12%  "
$ws.Range("D20").Value = "'7.34"
$ws.Range("E20").Value = "  +5.36%  "
$ws.Range("D21").Value = "'329.43"
$ws.Range("E21").Value = "  +1.72%  "
$ws.Range("D22").Value = "'4.22"
$ws.Range("E22").Value = "  +1.28%  "
$ws.Range("D23").Value = "'2.09"
$ws.Range("E23").Value = "  +15.04%  "
$ws.Range("D24").Value = "'1.00"
$ws.Range("E24").Value = "  +0.21%  "
$ws.Range("D25").Value = "'65.48"
$ws.Range("E25").Value = "  -2.30%  "
$ws.Range("D26").Value = "'616.40"
$ws.Range("E26").Value = "  +3.22%  "
$ws.Range("D27").Value = "'8.89"
$ws.Range("E27").Value = "  +3.68%  "
$ws.Range("E28").Value = "  +3.11%  "
$ws.Range("D29").Value = "2.580.73"
$ws.Range("E29").Value = "  +1.40%  "
$ws.Range("E30").Value = "  +5.02%  "
$ws.Range("D31").Value = "'0.998"
$ws.Range("E31").Value = "  -0.13%  "
$ws.Range("E32").Value = "  -2.25%  "
$ws.Range("D33").Value = "'1.91"
$ws.Range("E33").Value = "  +1.81%  "
$ws.Range("E34").Value = "  -2.15%  "
$ws.Range("D35").Value = "'5.21"
$ws.Range("E35").Value = "  +7.44%  "
$ws.Range("E36").Value = "  +1.96%  "
$ws.Range("E37").Value = "  +0.04%  "
$ws.Range("E38").Value = "  -0.40%  "
$ws.Range("D39").Value = "'18.88"
$ws.Range("E39").Value = "  +0.86%  "
$ws.Range("D40").Value = "'5.43"
$ws.Range("E40").Value = "  +1.49%  "
$ws.Range("D41").Value = "'147.17"
$ws.Range("E41").Value = "  -0.13%  "
$ws.Range("D42").Value = "'1.80"
$ws.Range("E42").Value = "  -0.91%  "
$ws.Range("D43").Value = "'2.62"
$ws.Range("E43").Value = "  +6.78%  "
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("E45").Value = "  +0.43%  "
$ws.Range("D46").Value = "'149.06"
$ws.Range("E46").Value = "  +0.68%  "
$ws.Range("D47").Value = "'3.78"
$ws.Range("E47").Value = "  +2.94%  "
$ws.Range("D48").Value = "'21.24"
$ws.Range("E48").Value = "  +3.61%  "
$ws.Range("E49").Value = "  -0.30%  "
$ws.Range("D50").Value = "'0.603"
$ws.Range("E50").Value = "  +0.28%  "
$ws.Range("E51").Value = "  +0.92%  "
